# "salvataggio query utente su file excel"
#
# The original template shipped with a single result sheet named
# "query3" holding the two-column header (Nome / Cognome) used by the
# user-query export. The save routine now emits the template with two
# worksheets - "Foglio1" and "Foglio2" - both carrying the same
# Nome/Cognome header row (and its styling), with "Foglio1" left as the
# active/selected sheet.

$wb = $excel.ActiveWorkbook

# The sheet that ships in the template today.
$source = $wb.Worksheets.Item(1)

# Build the two replacement sheets right after the source sheet so the
# final tab order is Foglio1, Foglio2.
$foglio1 = $wb.Worksheets.Add($null, $source)
$foglio1.Name = "Foglio1"

$foglio2 = $wb.Worksheets.Add($null, $foglio1)
$foglio2.Name = "Foglio2"

# Carry the header row (values + number/style formatting) over to both
# new sheets.
$source.Range("A1:B1").Copy($foglio1.Range("A1:B1"))
$source.Range("A1:B1").Copy($foglio2.Range("A1:B1"))

# Drop the now-superseded original sheet.
$source.Delete()

# Leave Foglio2's header selected too, then land back on Foglio1 as the
# active/selected sheet/range.
$f2 = $wb.Worksheets.Item("Foglio2")
$f2.Activate()
$f2.Range("A1:B1").Select()

$f1 = $wb.Worksheets.Item("Foglio1")
$f1.Activate()
$f1.Range("A1:B1").Select()
